# feat: add 2022-Q4 data
#
# 1. Insert a new "2022-Q4" worksheet (fund-holdings detail), positioned right
#    after "总计" and before "2022-Q3".
# 2. Insert a matching summary row for "2022-Q4" at the top of the "总计"
#    sheet's data table, pushing the older quarters down by one row.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")
$q3Sheet = $wb.Worksheets.Item("2022-Q3")

# ---------------------------------------------------------------------------
# 1. Create the "2022-Q4" sheet, placed immediately before "2022-Q3".
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($q3Sheet)
$newSheet.Name = "2022-Q4"

# Match page margins used by the rest of the workbook's sheets.
$newSheet.PageSetup.LeftMargin = 0.75 * 72
$newSheet.PageSetup.RightMargin = 0.75 * 72
$newSheet.PageSetup.TopMargin = 1 * 72
$newSheet.PageSetup.BottomMargin = 1 * 72
$newSheet.PageSetup.HeaderMargin = 0.5 * 72
$newSheet.PageSetup.FooterMargin = 0.5 * 72

# Header row (bold / bordered style, matching the other quarter sheets).
$newSheet.Cells.Item(1, 2).Value = "基金代码"
$newSheet.Cells.Item(1, 3).Value = "基金名称"
$newSheet.Cells.Item(1, 4).Value = "基金规模"
$newSheet.Cells.Item(1, 5).Value = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value = "仓位占比"
$newSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value = "仓位排名"
$q3Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# Fund rows. Columns B/D/E/F/G hold text (so values like "002810" or "3.84"
# keep their original formatting instead of being coerced into numbers),
# while A (index) and H (rank) are numeric.
$data = @(
    @("002810", "金信转型创新成长灵活配置混合", "3.84", "89.18", "5.79", "0.2223", 4),
    @("009490", "泰康科技创新一年定期开放混合", "2.21", "93.49", "2.56", "0.0566", 10),
    @("519097", "新华中小市值优选混合", "0.66", "70.51", "4.02", "0.0265", 4),
    @("519139", "海富通沪港深灵活配置混合", "0.67", "92.35", "2.88", "0.0193", 8)
)

$lastRow = $data.Count + 1
$newSheet.Range("B2:G$lastRow").NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rec = $data[$i]
    $newSheet.Cells.Item($row, 1).Value = $i
    $newSheet.Cells.Item($row, 2).Value = $rec[0]
    $newSheet.Cells.Item($row, 3).Value = $rec[1]
    $newSheet.Cells.Item($row, 4).Value = $rec[2]
    $newSheet.Cells.Item($row, 5).Value = $rec[3]
    $newSheet.Cells.Item($row, 6).Value = $rec[4]
    $newSheet.Cells.Item($row, 7).Value = $rec[5]
    $newSheet.Cells.Item($row, 8).Value = $rec[6]
}

# Column-A style (index column) matches the other quarter sheets.
$q3Sheet.Range("A2").Copy()
$newSheet.Range("A2:A$lastRow").PasteSpecial(-4122)
for ($i = 0; $i -lt $data.Count; $i++) {
    $newSheet.Cells.Item($i + 2, 1).Value = $i
}

# ---------------------------------------------------------------------------
# 2. Add the "2022-Q4" summary row to "总计", pushing the rest down.
# ---------------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("A3:D3").Copy()
$totalSheet.Range("A2:D2").PasteSpecial(-4122)

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q4"
$totalSheet.Cells.Item(2, 3).Value = 4
$totalSheet.Cells.Item(2, 4).Value = 0.32

# Renumber the 0-based index column for the rows that shifted down.
for ($r = 3; $r -le 6; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}
